# Update cryptos list values (prices and 1h volume %) to match the latest
# scrape. Also swaps RocketPoolETH / FraxShare rows (rank 46/47) with their
# updated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = '29.916.02'
$ws.Range("E2").Value = '  +0.80%  '

# Row 3
$ws.Range("D3").Value = '1.630.04'
$ws.Range("E3").Value = '  +1.69%  '

# Row 4
$ws.Range("E4").Value = '  +0.45%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '214.23'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

# Row 6
$ws.Range("E6").Value = '  -0.14%  '

# Row 7
$ws.Range("E7").Value = '  +0.46%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '28.25'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -3.54%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.257'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.50%  '

# Row 10
$ws.Range("E10").Value = '  +0.59%  '

# Row 11
$ws.Range("E11").Value = '  +0.21%  '

# Row 12
$ws.Range("D12").Value = '1.865.89'
$ws.Range("E12").Value = '  +1.72%  '

# Row 13
$ws.Range("D13").Value = '1.632.23'
$ws.Range("E13").Value = '  +1.68%  '

# Row 14
$ws.Range("E14").Value = '  +1.19%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '9.11'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +13.92%  '

# Row 16
$ws.Range("D16").Value = '29.927.75'
$ws.Range("E16").Value = '  +0.66%  '

# Row 17
$ws.Range("E17").Value = '  +1.32%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '63.92'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.26%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '241.25'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.51%  '

# Row 20
$ws.Range("E20").Value = '  +0.14%  '

# Row 21
$ws.Range("E21").Value = '  +0.37%  '

# Row 22
$ws.Range("E22").Value = '  +1.83%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.72'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.20%  '

# Row 24
$ws.Range("E24").Value = '  +2.34%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '160.98'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +3.51%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '15.48'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.21%  '

# Row 27
$ws.Range("E27").Value = '  +0.12%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '6.59'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.93%  '

# Row 29
$ws.Range("E29").Value = '  +0.42%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0485'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.25%  '

# Row 31
$ws.Range("E31").Value = '  +3.78%  '

# Row 32
$ws.Range("E32").Value = '  +3.55%  '

# Row 33
$ws.Range("E33").Value = '  -0.31%  '

# Row 34
$ws.Range("D34").Value = '1.423.56'
$ws.Range("E34").Value = '  -0.18%  '

# Row 35
$ws.Range("E35").Value = '  +4.17%  '

# Row 36
$ws.Range("E36").Value = '  -1.18%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.77'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.67%  '

# Row 38
$ws.Range("E38").Value = '  -0.08%  '

# Row 39
$ws.Range("E39").Value = '  -0.28%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '75.03'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +11.46%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.550'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.74%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.99'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.76%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.825'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.55%  '

# Row 44
$ws.Range("E44").Value = '  -0.47%  '

# Row 45
$ws.Range("E45").Value = '  +0.53%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.60%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '52.37'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -5.87%  '

# Row 48
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '5.35'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.54%  '

# Row 49
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '1.773.80'
$ws.Range("E49").Value = '  +1.82%  '

# Row 50
$ws.Range("E50").Value = '  +13.33%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '90.36'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +4.18%  '
